# Apply updated Betfair back/lay odds values for 2025-11-17 sheet (rows 2-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 6).Value = 2.86  # F2
$ws.Cells.Item(2, 9).Value = 3.25  # I2
$ws.Cells.Item(2, 11).Value = 3.15  # K2
$ws.Cells.Item(2, 12).Value = 1.57  # L2
$ws.Cells.Item(2, 14).Value = 2.66  # N2
$ws.Cells.Item(2, 15).Value = 1.52  # O2
$ws.Cells.Item(2, 16).Value = 1.55  # P2
$ws.Cells.Item(2, 17).Value = 2.6  # Q2
$ws.Cells.Item(2, 18).Value = 1.19  # R2
$ws.Cells.Item(2, 19).Value = 5.3  # S2
$ws.Cells.Item(2, 20).Value = 2.02  # T2
$ws.Cells.Item(2, 22).Value = 1.45  # V2
$ws.Cells.Item(2, 23).Value = 1.45  # W2
$ws.Cells.Item(2, 25).Value = 17.5  # Y2
$ws.Cells.Item(2, 28).Value = 1000  # AB2
$ws.Cells.Item(2, 31).Value = 130  # AE2
$ws.Cells.Item(2, 34).Value = 1000  # AH2
$ws.Cells.Item(2, 37).Value = 130  # AK2

# Row 3
$ws.Cells.Item(3, 6).Value = 2.02  # F3
$ws.Cells.Item(3, 7).Value = 2.24  # G3
$ws.Cells.Item(3, 8).Value = 4.2  # H3
$ws.Cells.Item(3, 9).Value = 5.2  # I3
$ws.Cells.Item(3, 10).Value = 2.96  # J3
$ws.Cells.Item(3, 11).Value = 3.2  # K3
$ws.Cells.Item(3, 12).Value = 1.56  # L3
$ws.Cells.Item(3, 13).Value = 1.11  # M3
$ws.Cells.Item(3, 14).Value = 2.66  # N3
$ws.Cells.Item(3, 15).Value = 1.51  # O3
$ws.Cells.Item(3, 16).Value = 1.56  # P3
$ws.Cells.Item(3, 17).Value = 2.5  # Q3
$ws.Cells.Item(3, 19).Value = 5.1  # S3
$ws.Cells.Item(3, 20).Value = 2.08  # T3
$ws.Cells.Item(3, 21).Value = 1.72  # U3
$ws.Cells.Item(3, 23).Value = 1.8  # W3
$ws.Cells.Item(3, 24).Value = 9.800000000000001  # X3
$ws.Cells.Item(3, 25).Value = 13  # Y3
$ws.Cells.Item(3, 26).Value = 55  # Z3
$ws.Cells.Item(3, 28).Value = 7.2  # AB3
$ws.Cells.Item(3, 29).Value = 8  # AC3
$ws.Cells.Item(3, 30).Value = 22  # AD3
$ws.Cells.Item(3, 32).Value = 12.5  # AF3
$ws.Cells.Item(3, 33).Value = 12  # AG3
$ws.Cells.Item(3, 36).Value = 80  # AJ3
$ws.Cells.Item(3, 37).Value = 46  # AK3
$ws.Cells.Item(3, 40).Value = 28  # AN3

# Row 4
$ws.Cells.Item(4, 6).Value = 4.4  # F4
$ws.Cells.Item(4, 7).Value = 5.1  # G4
$ws.Cells.Item(4, 8).Value = 2.08  # H4
$ws.Cells.Item(4, 9).Value = 2.22  # I4
$ws.Cells.Item(4, 10).Value = 2.96  # J4
$ws.Cells.Item(4, 11).Value = 3.25  # K4
$ws.Cells.Item(4, 12).Value = 1.66  # L4
$ws.Cells.Item(4, 13).Value = 1.14  # M4
$ws.Cells.Item(4, 15).Value = 1.63  # O4
$ws.Cells.Item(4, 16).Value = 1.44  # P4
$ws.Cells.Item(4, 17).Value = 3  # Q4
$ws.Cells.Item(4, 18).Value = 1.15  # R4
$ws.Cells.Item(4, 20).Value = 2.34  # T4
$ws.Cells.Item(4, 21).Value = 1.62  # U4
$ws.Cells.Item(4, 22).Value = 1.82  # V4
$ws.Cells.Item(4, 23).Value = 1.25  # W4
$ws.Cells.Item(4, 24).Value = 7.8  # X4
$ws.Cells.Item(4, 25).Value = 6.4  # Y4
$ws.Cells.Item(4, 26).Value = 11  # Z4
$ws.Cells.Item(4, 27).Value = 27  # AA4
$ws.Cells.Item(4, 28).Value = 11.5  # AB4
$ws.Cells.Item(4, 29).Value = 7.4  # AC4
$ws.Cells.Item(4, 30).Value = 12.5  # AD4
$ws.Cells.Item(4, 31).Value = 36  # AE4
$ws.Cells.Item(4, 32).Value = 32  # AF4
$ws.Cells.Item(4, 33).Value = 23  # AG4
$ws.Cells.Item(4, 34).Value = 42  # AH4
$ws.Cells.Item(4, 36).Value = 1000  # AJ4
$ws.Cells.Item(4, 37).Value = 120  # AK4
$ws.Cells.Item(4, 39).Value = 410  # AM4
$ws.Cells.Item(4, 41).Value = 36  # AO4

# Row 5
$ws.Cells.Item(5, 6).Value = 3.55  # F5
$ws.Cells.Item(5, 8).Value = 2.3  # H5
$ws.Cells.Item(5, 9).Value = 2.56  # I5
$ws.Cells.Item(5, 10).Value = 2.88  # J5
$ws.Cells.Item(5, 12).Value = 1.57  # L5
$ws.Cells.Item(5, 13).Value = 1.12  # M5
$ws.Cells.Item(5, 14).Value = 2.64  # N5
$ws.Cells.Item(5, 15).Value = 1.51  # O5
$ws.Cells.Item(5, 16).Value = 1.56  # P5
$ws.Cells.Item(5, 17).Value = 2.58  # Q5
$ws.Cells.Item(5, 18).Value = 1.21  # R5
$ws.Cells.Item(5, 21).Value = 1.81  # U5
$ws.Cells.Item(5, 22).Value = 1.65  # V5

# Row 6
$ws.Cells.Item(6, 6).Value = 2.64  # F6
$ws.Cells.Item(6, 8).Value = 3  # H6
$ws.Cells.Item(6, 9).Value = 3.45  # I6
$ws.Cells.Item(6, 11).Value = 3.2  # K6
$ws.Cells.Item(6, 12).Value = 1.63  # L6
$ws.Cells.Item(6, 13).Value = 1.14  # M6
$ws.Cells.Item(6, 14).Value = 2.44  # N6
$ws.Cells.Item(6, 17).Value = 2.74  # Q6
$ws.Cells.Item(6, 19).Value = 5.8  # S6
$ws.Cells.Item(6, 21).Value = 1.71  # U6
$ws.Cells.Item(6, 22).Value = 1.41  # V6
$ws.Cells.Item(6, 23).Value = 1.5  # W6
$ws.Cells.Item(6, 29).Value = 8.199999999999999  # AC6
$ws.Cells.Item(6, 30).Value = 1000  # AD6

# Row 7
$ws.Cells.Item(7, 6).Value = 2.54  # F7
$ws.Cells.Item(7, 7).Value = 2.62  # G7
$ws.Cells.Item(7, 8).Value = 3.25  # H7
$ws.Cells.Item(7, 9).Value = 3.4  # I7
$ws.Cells.Item(7, 10).Value = 3.2  # J7
$ws.Cells.Item(7, 12).Value = 1.52  # L7
$ws.Cells.Item(7, 14).Value = 3.05  # N7
$ws.Cells.Item(7, 16).Value = 1.66  # P7
$ws.Cells.Item(7, 17).Value = 2.4  # Q7
$ws.Cells.Item(7, 18).Value = 1.26  # R7
$ws.Cells.Item(7, 20).Value = 2  # T7
$ws.Cells.Item(7, 21).Value = 1.93  # U7
$ws.Cells.Item(7, 22).Value = 1.41  # V7
$ws.Cells.Item(7, 23).Value = 1.61  # W7
$ws.Cells.Item(7, 24).Value = 10  # X7
$ws.Cells.Item(7, 25).Value = 11  # Y7
$ws.Cells.Item(7, 26).Value = 22  # Z7
$ws.Cells.Item(7, 29).Value = 7  # AC7
$ws.Cells.Item(7, 30).Value = 14.5  # AD7
$ws.Cells.Item(7, 31).Value = 46  # AE7
$ws.Cells.Item(7, 33).Value = 12  # AG7
$ws.Cells.Item(7, 35).Value = 70  # AI7
$ws.Cells.Item(7, 37).Value = 34  # AK7
$ws.Cells.Item(7, 39).Value = 150  # AM7
$ws.Cells.Item(7, 40).Value = 34  # AN7
$ws.Cells.Item(7, 41).Value = 55  # AO7

# Row 8
$ws.Cells.Item(8, 12).Value = 1.81  # L8
$ws.Cells.Item(8, 13).Value = 1.2  # M8
$ws.Cells.Item(8, 14).Value = 2.12  # N8
$ws.Cells.Item(8, 15).Value = 1.83  # O8
$ws.Cells.Item(8, 16).Value = 1.34  # P8
$ws.Cells.Item(8, 17).Value = 3.7  # Q8
$ws.Cells.Item(8, 20).Value = 2.58  # T8
$ws.Cells.Item(8, 21).Value = 1.54  # U8
$ws.Cells.Item(8, 22).Value = 1.47  # V8
$ws.Cells.Item(8, 23).Value = 1.44  # W8
$ws.Cells.Item(8, 24).Value = 5.9  # X8
$ws.Cells.Item(8, 25).Value = 7.4  # Y8
$ws.Cells.Item(8, 26).Value = 17  # Z8
$ws.Cells.Item(8, 27).Value = 65  # AA8
$ws.Cells.Item(8, 28).Value = 7.2  # AB8
$ws.Cells.Item(8, 29).Value = 7  # AC8
$ws.Cells.Item(8, 30).Value = 16.5  # AD8
$ws.Cells.Item(8, 32).Value = 18  # AF8
$ws.Cells.Item(8, 34).Value = 36  # AH8
$ws.Cells.Item(8, 35).Value = 130  # AI8
$ws.Cells.Item(8, 36).Value = 1000  # AJ8
$ws.Cells.Item(8, 37).Value = 70  # AK8
$ws.Cells.Item(8, 38).Value = 140  # AL8
$ws.Cells.Item(8, 39).Value = 1000  # AM8
$ws.Cells.Item(8, 40).Value = 110  # AN8
$ws.Cells.Item(8, 41).Value = 100  # AO8

# Row 9
$ws.Cells.Item(9, 6).Value = 2.2  # F9
$ws.Cells.Item(9, 7).Value = 2.28  # G9
$ws.Cells.Item(9, 8).Value = 4  # H9
$ws.Cells.Item(9, 10).Value = 3.15  # J9
$ws.Cells.Item(9, 12).Value = 1.61  # L9
$ws.Cells.Item(9, 14).Value = 2.62  # N9
$ws.Cells.Item(9, 15).Value = 1.58  # O9
$ws.Cells.Item(9, 17).Value = 2.78  # Q9
$ws.Cells.Item(9, 19).Value = 5.9  # S9
$ws.Cells.Item(9, 20).Value = 2.22  # T9
$ws.Cells.Item(9, 21).Value = 1.73  # U9
$ws.Cells.Item(9, 22).Value = 1.31  # V9
$ws.Cells.Item(9, 23).Value = 1.78  # W9
$ws.Cells.Item(9, 24).Value = 8.6  # X9
$ws.Cells.Item(9, 25).Value = 12.5  # Y9
$ws.Cells.Item(9, 26).Value = 38  # Z9
$ws.Cells.Item(9, 27).Value = 130  # AA9
$ws.Cells.Item(9, 29).Value = 7.8  # AC9
$ws.Cells.Item(9, 30).Value = 25  # AD9
$ws.Cells.Item(9, 31).Value = 1000  # AE9
$ws.Cells.Item(9, 32).Value = 13  # AF9
$ws.Cells.Item(9, 33).Value = 12  # AG9
$ws.Cells.Item(9, 35).Value = 120  # AI9
$ws.Cells.Item(9, 38).Value = 1000  # AL9
$ws.Cells.Item(9, 41).Value = 140  # AO9

# Row 10
$ws.Cells.Item(10, 6).Value = 2.08  # F10
$ws.Cells.Item(10, 7).Value = 2.16  # G10
$ws.Cells.Item(10, 8).Value = 4.5  # H10
$ws.Cells.Item(10, 9).Value = 4.9  # I10
$ws.Cells.Item(10, 10).Value = 3.15  # J10
$ws.Cells.Item(10, 12).Value = 1.58  # L10
$ws.Cells.Item(10, 14).Value = 2.74  # N10
$ws.Cells.Item(10, 16).Value = 1.55  # P10
$ws.Cells.Item(10, 17).Value = 2.68  # Q10
$ws.Cells.Item(10, 20).Value = 2.2  # T10
$ws.Cells.Item(10, 21).Value = 1.74  # U10
$ws.Cells.Item(10, 23).Value = 1.87  # W10
$ws.Cells.Item(10, 24).Value = 9.199999999999999  # X10
$ws.Cells.Item(10, 25).Value = 14  # Y10
$ws.Cells.Item(10, 28).Value = 7.4  # AB10
$ws.Cells.Item(10, 30).Value = 24  # AD10
$ws.Cells.Item(10, 33).Value = 12  # AG10
$ws.Cells.Item(10, 34).Value = 40  # AH10
$ws.Cells.Item(10, 36).Value = 980  # AJ10
$ws.Cells.Item(10, 37).Value = 34  # AK10
$ws.Cells.Item(10, 39).Value = 230  # AM10
$ws.Cells.Item(10, 40).Value = 38  # AN10

# Row 11
$ws.Cells.Item(11, 6).Value = 2.48  # F11
$ws.Cells.Item(11, 7).Value = 2.52  # G11
$ws.Cells.Item(11, 9).Value = 3.45  # I11
$ws.Cells.Item(11, 10).Value = 3.25  # J11
$ws.Cells.Item(11, 11).Value = 3.3  # K11
$ws.Cells.Item(11, 12).Value = 1.64  # L11
$ws.Cells.Item(11, 14).Value = 2.6  # N11
$ws.Cells.Item(11, 15).Value = 1.59  # O11
$ws.Cells.Item(11, 16).Value = 1.53  # P11
$ws.Cells.Item(11, 17).Value = 2.76  # Q11
$ws.Cells.Item(11, 18).Value = 1.18  # R11
$ws.Cells.Item(11, 19).Value = 5.9  # S11
$ws.Cells.Item(11, 20).Value = 2.22  # T11
$ws.Cells.Item(11, 21).Value = 1.74  # U11
$ws.Cells.Item(11, 22).Value = 1.41  # V11
$ws.Cells.Item(11, 23).Value = 1.66  # W11
$ws.Cells.Item(11, 24).Value = 8.6  # X11
$ws.Cells.Item(11, 25).Value = 9.4  # Y11
$ws.Cells.Item(11, 28).Value = 7.4  # AB11
$ws.Cells.Item(11, 29).Value = 7.6  # AC11
$ws.Cells.Item(11, 30).Value = 15.5  # AD11
$ws.Cells.Item(11, 31).Value = 250  # AE11
$ws.Cells.Item(11, 32).Value = 17.5  # AF11
$ws.Cells.Item(11, 34).Value = 26  # AH11
$ws.Cells.Item(11, 35).Value = 85  # AI11
$ws.Cells.Item(11, 36).Value = 40  # AJ11
$ws.Cells.Item(11, 37).Value = 38  # AK11
$ws.Cells.Item(11, 39).Value = 200  # AM11
$ws.Cells.Item(11, 40).Value = 42  # AN11
$ws.Cells.Item(11, 41).Value = 85  # AO11
